# "secondbranch": drop the opening sentence "This is demo2 file for
# practice." (originally split across three separate runs) and keep only
# the trailing sentence, as a single run of replacement text.
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "This is demo2 file for practice. We are creating the pull requests.",
    $true,                               # MatchCase
    $false,                              # MatchWholeWord
    $false,                              # MatchWildcards
    $false,                              # MatchSoundsLike
    $false,                              # MatchAllWordForms
    $true,                               # Forward
    1,                                   # Wrap (wdFindContinue)
    $false,                              # Format
    "We are creating the pull requests.", # ReplaceWith
    2                                    # Replace (wdReplaceAll)
)
